# Add "d_name" / "mu_name" coefficient-lookup columns (G:H) to each of the
# four parameter-lookup sheets, filling in the mu_name values that line up
# with the existing mu_num column, and updating each non-active sheet's
# selection to the newly added mu_name range.

$wb = $excel.ActiveWorkbook

# --- weibull (mu_num goes 1,2,3,-,4,5 on rows 2-7) ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("G1").Value = "d_name"
$ws.Range("H1").Value = "mu_name"
$ws.Range("H2").Value = "mu_1"
$ws.Range("H3").Value = "mu_2"
$ws.Range("H4").Value = "mu_4"
$ws.Range("H6").Value = "mu_5"
$ws.Range("H7").Value = "mu_6"

# --- gompertz (same row/mu_num layout as weibull) ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("G1").Value = "d_name"
$ws.Range("H1").Value = "mu_name"
$ws.Range("H2").Value = "mu_1"
$ws.Range("H3").Value = "mu_2"
$ws.Range("H4").Value = "mu_4"
$ws.Range("H6").Value = "mu_5"
$ws.Range("H7").Value = "mu_6"
[void]$ws.Range("H2:H7").Select()

# --- fracpoly1 (mu_num goes 1,2,3,4,-,-,5,6,- on rows 2-10) ---
$ws = $wb.Worksheets.Item("fracpoly1")
$ws.Range("G1").Value = "d_name"
$ws.Range("H1").Value = "mu_name"
$ws.Range("H2").Value = "mu_1"
$ws.Range("H3").Value = "mu_2"
$ws.Range("H4").Value = "mu_3"
$ws.Range("H5").Value = "mu_4"
$ws.Range("H8").Value = "mu_5"
$ws.Range("H9").Value = "mu_6"
[void]$ws.Range("H2:H9").Select()

# --- fracpoly2 (same row/mu_num layout as fracpoly1) ---
$ws = $wb.Worksheets.Item("fracpoly2")
$ws.Range("G1").Value = "d_name"
$ws.Range("H1").Value = "mu_name"
$ws.Range("H2").Value = "mu_1"
$ws.Range("H3").Value = "mu_2"
$ws.Range("H4").Value = "mu_3"
$ws.Range("H5").Value = "mu_4"
$ws.Range("H8").Value = "mu_5"
$ws.Range("H9").Value = "mu_6"
[void]$ws.Range("H2:H9").Select()

# Restore the original active sheet/tab selection (weibull).
$ws1 = $wb.Worksheets.Item("weibull")
$ws1.Activate()
[void]$ws1.Range("A2").Select()
